# Sprint Plan - add new use cases (Passwort fuer Zielsetzung / Ansichtwechsel
# Admin App / Passwort fuer Zuruecksetzen) under the 6th sprint's "User
# Stories" column (J:K), and highlight the now-completed / in-focus stories
# with the new green fill.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. New rows for the two extra use-cases that didn't have template rows
#    yet (rows 6,7,8 - sheet used to stop at row 5).
# ---------------------------------------------------------------------

# Row 6 carries the 3rd new use case and needs the same look (thin border
# all around + centered) as the rest of the grid - clone it from an
# existing "blank" template cell before putting content in it.
$ws.Range("H4").Copy() | Out-Null
$ws.Range("J6").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("K6").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$ws.Range("H4").Copy() | Out-Null
$ws.Range("J7").PasteSpecial(-4122) | Out-Null
$ws.Range("K7").PasteSpecial(-4122) | Out-Null
$ws.Range("J8").PasteSpecial(-4122) | Out-Null
$ws.Range("K8").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$ws.Range("J6:K6").Merge()
$ws.Range("J7:K7").Merge()
$ws.Range("J8:K8").Merge()

# ---------------------------------------------------------------------
# 2. Fill in the new use-case text.
# ---------------------------------------------------------------------
$ws.Range("J4").Value = "Passwort für Zielsetzung"
$ws.Range("J5").Value = "Ansichtwechsel Admin App"
$ws.Range("J6").Value = "Passwort für Zurücksetzen"

# ---------------------------------------------------------------------
# 3. Highlight the relevant user-story boxes with the new green fill
#    (RGB 0,176,80 -> hex 00B050) used to flag these use cases.
# ---------------------------------------------------------------------
$green = 0 + (0xB0 * 256) + (0x50 * 65536)

$ws.Range("D3:E3").Interior.Color = $green
$ws.Range("F3:G3").Interior.Color = $green
$ws.Range("H3:I3").Interior.Color = $green
$ws.Range("F4:G4").Interior.Color = $green
$ws.Range("J4:K4").Interior.Color = $green
$ws.Range("J5:K5").Interior.Color = $green
$ws.Range("J6:K6").Interior.Color = $green

# ---------------------------------------------------------------------
# 4. Update the visible selection to reflect where the edit happened.
# ---------------------------------------------------------------------
$ws.Range("L7").Select() | Out-Null
